$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12-74 shift down to 13-75.
$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with the new weekly price record.
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 44819
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 100112013
$ws.Cells.Item(12, 7).Value = "Alcachofa"
$ws.Cells.Item(12, 8).Value = "Argentina(o)"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 110
$ws.Cells.Item(12, 11).Value = 12000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 13364
$ws.Cells.Item(12, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 267
$ws.Cells.Item(12, 17).Value = 50
$ws.Cells.Item(12, 18).Value = "Hortaliza"
